# 7.1.2.xlsx edit: add 2020 column (column E) data to the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlRight = -4152
$xlCenter = -4108

function Set-EFormat([string]$row) {
    # Copy D column formatting onto the matching E cell, then re-apply the
    # numeric format + alignment that column E (2020 data) uses.
    $ws.Range("D$row").Copy()
    $ws.Range("E$row").PasteSpecial($xlPasteFormats)
    $ws.Range("E$row").NumberFormat = "0.0"
    $ws.Range("E$row").HorizontalAlignment = $xlRight
    $ws.Range("E$row").VerticalAlignment = $xlCenter
}

# Row 3: header year value, re-use D3's exact style (thick-bottom border row)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Value = 2020

# Row 4: top-level total row (has top border)
Set-EFormat "4"
$ws.Range("E4").Value = 22.1

# Row 5: Urbanisation header row (stays empty)
Set-EFormat "5"

# Row 6-7: urban / rural
Set-EFormat "6"
$ws.Range("E6").Value = 52.7
Set-EFormat "7"
$ws.Range("E7").Value = 5

# Row 8: Education header row (stays empty)
Set-EFormat "8"

# Rows 9-15: education levels
Set-EFormat "9"
$ws.Range("E9").Value = 4.8
Set-EFormat "10"
$ws.Range("E10").Value = 15.8
Set-EFormat "11"
$ws.Range("E11").Value = 13.5
Set-EFormat "12"
$ws.Range("E12").Value = 9.6
Set-EFormat "13"
$ws.Range("E13").Value = 2.7
Set-EFormat "14"
$ws.Range("E14").Value = 14.7
Set-EFormat "15"
$ws.Range("E15").Value = 18.2

# Rows 16-17: city totals
Set-EFormat "16"
$ws.Range("E16").Value = 74
Set-EFormat "17"
$ws.Range("E17").Value = 35.1

# Row 18: By territory header row (stays empty)
Set-EFormat "18"

# Rows 19-23: regions with no 2020 data ("-")
Set-EFormat "19"
$ws.Range("E19").Value = "-"
Set-EFormat "20"
$ws.Range("E20").Value = "-"
Set-EFormat "21"
$ws.Range("E21").Value = "-"
Set-EFormat "22"
$ws.Range("E22").Value = "-"
Set-EFormat "23"
$ws.Range("E23").Value = "-"

# Row 24: Wealth index quintile header row (stays empty)
Set-EFormat "24"

# Rows 25-29: wealth quintiles with no 2020 data ("-")
Set-EFormat "25"
$ws.Range("E25").Value = "-"
Set-EFormat "26"
$ws.Range("E26").Value = "-"
Set-EFormat "27"
$ws.Range("E27").Value = "-"
Set-EFormat "28"
$ws.Range("E28").Value = "-"
Set-EFormat "29"
$ws.Range("E29").Value = "-"

# Update the active selection shown in the sheet view
$ws.Range("J24").Select()
